$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 44476
$ws.Range("D8").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 100112042
$ws.Range("G8").Value = "Locoto"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 2200
$ws.Range("L8").Value = 2200
$ws.Range("M8").Value = 2200
$ws.Range("N8").Value = '$/kilo'
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 2200
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
